# Update the workbook to reflect running the logbook processing "at CC"
# (CC = Computing Center), per commit message "update to use at CC".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The "path" column (D) held a local macOS path that pointed at the
#    author's own machine. Every data row (D2:D125) shares this same
#    string value, so replacing the text updates them all at once.
$oldPath = "/Users/dagoret/DATA/AuxTelData2021/holo/quickLookExp/2021-07-07"
$newPath = "/sps/lsst/groups/auxtel/data/2021/holo/quickLookExp/2021-07-07"

$ws.Cells.Replace($oldPath, $newPath)

# 2) Row 98 (data row "96") was missing its Obj-posXpix / Obj-posYpix /
#    run / quality values that are present on the neighboring rows.
#    Fill them in to match that pattern.
$ws.Range("R98").Value = 300
$ws.Range("S98").Value = 1600
$ws.Range("W98").Value = 1
$ws.Range("X98").Value = 1
